$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.225.61"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.600.70"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.18"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.484"
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0613"
$ws.Range("E9").Value = "  -0.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.11"
$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +2.10%  "

$ws.Range("D12").Value = "1.823.70"
$ws.Range("E12").Value = "  -0.26%  "

$ws.Range("D13").Value = "1.602.94"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  +0.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.513"
$ws.Range("E15").Value = "  +0.80%  "

$ws.Range("D16").Value = "26.212.91"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.22"
$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.21"
$ws.Range("E20").Value = "  +1.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.00"
$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.92"
$ws.Range("E24").Value = "  +9.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.61"
$ws.Range("E25").Value = "  +1.86%  "

$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").Value = "  -7.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.16"
$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("E29").Value = "  +1.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0490"
$ws.Range("E30").Value = "  +3.53%  "

$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.15"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("E33").Value = "  -3.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.42"
$ws.Range("E34").Value = "  +2.77%  "

$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("D36").Value = "1.153.10"
$ws.Range("E36").Value = "  +4.10%  "

$ws.Range("E37").Value = "  +7.82%  "

$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("E39").Value = "  -0.68%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.789"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.497"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("D44").Value = "1.739.24"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.91"
$ws.Range("E45").Value = "  -1.07%  "

$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.02"
$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "0.0₇0951"
$ws.Range("E51").Value = "  -7.47%  "
